$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the municipio-nombre column (D) mapping metadata to the new curated dimension
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Update the aragon column (E) mapping metadata to the new curated dimension
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Comunidad"

# Remove the old row referencing the mapping-aragon.xlsx file (no longer needed)
$ws.Range("E5").EntireRow.Delete()
